# Correção do salvamento dos relatorios.
# Reorders the header columns (Valor Total Pago / Filme / Cliente) and
# (re)populates the report rows that the faulty save logic had dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Codigo Reserva"
$ws.Range("B1").Value = "Data de Entrada do Valor"
$ws.Range("C1").Value = "Valor Total Pago"
$ws.Range("D1").Value = "Filme"
$ws.Range("E1").Value = "Cliente"

# --- Data rows ----------------------------------------------------------
$data = @(
    @(9,  "2021-10-22T12:02:18-03:00", 3.75,  "Homem de Ferro 1",          "Jose Gabriel Batista Machado"),
    @(10, "2021-10-22T12:02:33-03:00", 3.6,   "Star Wars: The Clone Wars", "Murilo jose"),
    @(12, "2021-10-22T14:17:36-03:00", 3.6,   "O incrível Hulk",           "Murilo jose"),
    @(13, "2021-10-22T14:40:58-03:00", 3.75,  "Homem de Ferro 1",          "Jose Gabriel Batista Machado"),
    @(15, "2021-10-22T16:26:04-03:00", 2.85,  "Rei Arthur",                "Isaque Almeida Damascen"),
    @(17, "2021-10-25T14:25:27-03:00", 25.35, "Homem de Ferro 1",          "Jose Gabriel Batista Machado"),
    @(8,  "2021-10-26T10:45:31-03:00", 2.55,  "A Culpa e das Estrelas",    "Flavia Cristina nascimento"),
    @(7,  "2021-10-26T16:54:24-03:00", 8.6,   "A Culpa e das Estrelas",    "Flavia Cristina nascimento"),
    @(21, "2021-10-29T13:09:11-03:00", 0.15,  "O incrível Hulk",           "Jose Gabriel Batista Machado"),
    @(5,  "2021-10-29T13:11:13-03:00", 37.15, "Homem de Ferro 1",          "Jose Gabriel Batista Machado"),
    @(6,  "2021-10-29T13:11:18-03:00", 36.8,  "Star Wars: The Clone Wars", "Murilo jose")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
